$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowIndex = 26

$values = @(
    "622229",
    "Facilitator",
    "kyla",
    "pineda",
    "asdhads@gmail.com",
    "09182388232",
    "filipino",
    "Catholic",
    "Male",
    "Single",
    "44",
    "no",
    "kyadl",
    "sadasd"
)

for ($col = 1; $col -le $values.Length; $col++) {
    $cell = $ws.Cells.Item($rowIndex, $col)
    # Leading apostrophe forces the value to be stored as text so that
    # numeric-looking strings (IDs, phone numbers, ages) keep exact
    # formatting (e.g. preserved leading zeros) instead of becoming numbers.
    $cell.Formula = "'" + $values[$col - 1]
    # Reset the cell style back to Normal/default so no extra number
    # formatting (e.g. the Text "@" format) is left applied to the cell,
    # matching the unstyled data rows already in the sheet.
    $cell.Style = "Normal"
}
